# 7.1.2.xlsx — add the "2020" column (E) of data next to the existing
# "2018" column (D), mirroring the number-formatting/border treatment
# already used on column D, then move the active selection to J24
# (as recorded by the authoring Excel session).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122
# xlRight
$xlRight = -4152
# xlCenter
$xlCenter = -4108

# --- Row 3: header year "2020", reuse the exact header style from D3 ---
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial($xlPasteFormats)
$ws.Range("E3").Value = 2020

# --- Row 4: first data row, border/format matches column D's top row ---
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial($xlPasteFormats)
$ws.Range("E4").NumberFormat = "0.0"
$ws.Range("E4").HorizontalAlignment = $xlRight
$ws.Range("E4").VerticalAlignment = $xlCenter
$ws.Range("E4").Value = 22.1

# --- Rows 5-18: regular data rows (same base style as column D) ---
$data5_18 = @{
    5  = $null
    6  = 52.7
    7  = 5
    8  = $null
    9  = 4.8
    10 = 15.8
    11 = 13.5
    12 = 9.6
    13 = 2.7
    14 = 14.7
    15 = 18.2
    16 = 74
    17 = 35.1
    18 = $null
}
foreach ($r in 5..18) {
    $ws.Range("D5").Copy()
    $ws.Range("E$r").PasteSpecial($xlPasteFormats)
    $ws.Range("E$r").NumberFormat = "0.0"
    $ws.Range("E$r").HorizontalAlignment = $xlRight
    $ws.Range("E$r").VerticalAlignment = $xlCenter
    $v = $data5_18[$r]
    if ($v -ne $null) {
        $ws.Range("E$r").Value = $v
    }
}

# --- Rows 19-23 and 25-28: "no data" rows, shown as "-" ---
foreach ($r in (19..23) + (25..28)) {
    $ws.Range("D26").Copy()
    $ws.Range("E$r").PasteSpecial($xlPasteFormats)
    $ws.Range("E$r").NumberFormat = "0.0"
    $ws.Range("E$r").HorizontalAlignment = $xlRight
    $ws.Range("E$r").VerticalAlignment = $xlCenter
    $ws.Range("E$r").Value = "-"
}

# --- Row 24: section header row, stays empty but gains the new format ---
$ws.Range("D24").Copy()
$ws.Range("E24").PasteSpecial($xlPasteFormats)
$ws.Range("E24").NumberFormat = "0.0"
$ws.Range("E24").HorizontalAlignment = $xlRight
$ws.Range("E24").VerticalAlignment = $xlCenter

# --- Row 29: bottom-border total row, shown as "-" ---
$ws.Range("D29").Copy()
$ws.Range("E29").PasteSpecial($xlPasteFormats)
$ws.Range("E29").NumberFormat = "0.0"
$ws.Range("E29").HorizontalAlignment = $xlRight
$ws.Range("E29").VerticalAlignment = $xlCenter
$ws.Range("E29").Value = "-"

# --- restore the authoring session's active selection ---
$ws.Range("J24").Select()
